$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# Update existing descriptions (column G = Description) to match new wording
$ws.Range("G6").Value = "J'ai créé les différentes features du projet"
$ws.Range("G7").Value = "J'ai ajouté quelques tâches au premier sprint"
$ws.Range("G8").Value = "J'ai ajouté quelques tâches au premier sprint"
$ws.Range("G9").Value = "J'ai finalizé le MCD"

# Copy the formatting of the last existing table row so new rows inherit
# the same cell styles (borders, number formats, etc.)
$lastDataRow = $lo.ListRows.Item($lo.ListRows.Count).Range
$lastDataRow.Copy()

# --- New row 10 ---
$row10 = $lo.ListRows.Add()
$row10.Range.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$rng10 = $row10.Range
$rng10.Item(1,1).Value = 45050
$rng10.Item(1,2).Value = 0.56527777777777777
$rng10.Item(1,3).Value = 0.62847222222222221
$rng10.Item(1,4).Formula = "=Tableau4[[#This Row],[Heure fin]]-Tableau4[[#This Row],[Heure début]]"
$rng10.Item(1,5).Value = "CPNV"
$rng10.Item(1,6).Value = "MLD"
$rng10.Item(1,7).Value = "J'ai commencé le MLD"

# --- New row 11 ---
$lastDataRow.Copy()
$row11 = $lo.ListRows.Add()
$row11.Range.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$rng11 = $row11.Range
$rng11.Item(1,1).Value = 45050
$rng11.Item(1,2).Value = 0.64097222222222217
$rng11.Item(1,3).Value = 0.66597222222222219
$rng11.Item(1,4).Formula = "=Tableau4[[#This Row],[Heure fin]]-Tableau4[[#This Row],[Heure début]]"
$rng11.Item(1,5).Value = "CPNV"
$rng11.Item(1,6).Value = "MLD"
$rng11.Item(1,7).Value = "J'ai continué le MLD"

$ws.Range("C11").Select()

$wb.Save()
